$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before A. This shifts the existing columns A:K to
#    B:L (values, styles AND merged-cell ranges all shift automatically).
# ---------------------------------------------------------------------------
$ws.Columns("A").Insert()

# ---------------------------------------------------------------------------
# 2. Reorder the three "Libro" blocks (rows 2-7, 8-13, 14-19) so they are
#    sorted descending by "Numero de Valoraciones":
#       new block1 (rows 2-7)   <- old block3 (rows 14-19)  Libro 42  (156)
#       new block2 (rows 8-13)  <- old block1 (rows 2-7)    Libro 18  (132)
#       new block3 (rows 14-19) <- old block2 (rows 8-13)   Libro 25  (68)
#
#    Columns H:L ("Opinion del 1 al 6", fi, hi, Fi, Hi) are NOT merged, so a
#    plain range Copy/PasteSpecial(values) can reorder them safely (this does
#    not disturb any merged ranges, which all live in columns A:G).
#    A scratch area far below the data (row 200+) is used to stage the
#    three-way rotation so no block is overwritten before it's been copied.
# ---------------------------------------------------------------------------
$ws.Range("H14:L19").Copy()
$ws.Range("H200").PasteSpecial(-4163)
$ws.Range("H2:L7").Copy()
$ws.Range("H206").PasteSpecial(-4163)
$ws.Range("H8:L13").Copy()
$ws.Range("H212").PasteSpecial(-4163)

$ws.Range("H200:L205").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("H206:L211").Copy()
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("H212:L217").Copy()
$ws.Range("H14").PasteSpecial(-4163)

$ws.Range("H200:L217").Clear()

# ---------------------------------------------------------------------------
# 3. Columns B:G (Libro, Numero de Valoraciones, Media, Cuasidesviacion,
#    Mediana, Moda) are each merged in three 6-row blocks (B2:B7, B8:B13,
#    B14:B19, ...). Writing directly to the merge-anchor cell's .Value does
#    NOT break the merge (unlike Copy/PasteSpecial on a merged destination),
#    so the block reorder is done with literal writes to the anchor cells.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Libro 42 "
$ws.Range("C2").Value = 156
$ws.Range("D2").Value = 5.564102564102564
$ws.Range("E2").Value = 0.7802355903888978
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 6

$ws.Range("B8").Value = "Libro 18"
$ws.Range("C8").Value = 132
$ws.Range("D8").Value = 5.727272727272728
$ws.Range("E8").Value = 0.6189243203857999
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 6

$ws.Range("B14").Value = "Libro 25"
$ws.Range("C14").Value = 68
$ws.Range("D14").Value = 5.823529411764706
$ws.Range("E14").Value = 0.5166244188642394
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 6

# ---------------------------------------------------------------------------
# 4. Add the new "Especialidad" column in A: a header plus a single value
#    ("Economia") merged/spanning all 18 data rows (A2:A19).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Especialidad"
$ws.Range("A2:A19").Merge()
$ws.Range("A2").Value = "Economía"

# Normalize formatting: re-apply the shared header/body style (already used
# uniformly by columns B:G) onto column A so every cell - including the
# interior of the new merge - carries the same style instead of the
# per-edge borders Excel synthesizes automatically on Merge().
$ws.Range("B1:B19").Copy()
$ws.Range("A1:A19").PasteSpecial(-4122)
